$wb = $excel.ActiveWorkbook

# Sheets "展览" and "全部类型" both contain the same data table and both
# need the "想去人数" (F column) values refreshed to the newer scrape.
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 1939
    $ws.Range("F3").Value = 360
    $ws.Range("F4").Value = 1199
    $ws.Range("F5").Value = 1340
    $ws.Range("F7").Value = 6053
}
